# Add data for 2022-02-28 (carjacking-by-neighborhood-by-month)
# Updates the "through February 19" running total column (column B) to
# "through February 20" and bumps the affected neighborhood counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and refresh the running-total header text.
$ws.Name = "Through 2022-02-20"
$ws.Range("B1").Value = "February 2022 (through February 20)"

# Austin
$ws.Range("D3").Value = 9
$ws.Range("F3").Value = 5

# Auburn Gresham
$ws.Range("B7").Value = 3
$ws.Range("P7").Value = 1

# North Lawndale
$ws.Range("B8").Value = 6

# Chicago Lawn
$ws.Range("F16").Value = 1

# Near South Side
$ws.Range("B19").Value = 2

# Morgan Park
$ws.Range("B25").Value = 1

# Roseland
$ws.Range("B27").Value = 1

# West Loop
$ws.Range("B34").Value = 4

# Albany Park
$ws.Range("B40").Value = ""
$ws.Range("H40").Value = 1

# River North
$ws.Range("F42").Value = 1

# Little Village
$ws.Range("N47").Value = 2

# Avalon Park
$ws.Range("D55").Value = 2
